$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, identical change applied to both
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.
$updates = @{
    2  = 233
    3  = 261
    4  = 275
    5  = 814
    6  = 263
    7  = 6460
    9  = 71
    11 = 73
    12 = 33
    14 = 2
    15 = 203
    16 = 512
    17 = 50
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
